$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neo4jModel")

# Remove the worksheet protection so the template cells can be edited.
# (Protect()+Unprotect() is needed to fully clear the stored protection
# record, since a bare Unprotect() on an already-protected sheet is a
# no-op here.)
$ws.Protect()
$ws.Unprotect()

# Rename the "Relations" table title to reflect it now covers nodes too.
$ws.Range("A1").Value = "Table 1. Nodes and Relations"

# Clear out the sample "Relations" table rows (kept headers in rows 1-5,
# but the example data rows 6-9 are wiped back to blank cells).
$ws.Range("A6:C9").ClearContents()

# Row 9 no longer needs the taller wrap-text row height that the removed
# long "protocol" sample text required, so let it fall back to the
# worksheet's default row height.
$ws.Rows("9:9").AutoFit()

# Clear out the sample "Node P:V Pairs" table rows 8-12 similarly.
$ws.Range("E8:G12").ClearContents()

# Move the active selection to B9, matching where the editor left off.
$ws.Range("B9").Select()
